# "added 1 point in exce Docemnt in future power"
#
# The sheet previously ended with:
#   A5 = "note: 1 Delivery Note or Receipt Note will be Controlled in sales and purchase"
#   B5 = 3000
#
# A new bullet point is inserted above the closing note, the trailing
# numeric value is moved down and bumped, and a "Cost" header is added
# next to the existing bold "Customization" header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the old closing note down to row 6, and put the new bullet in
# its place at row 5 (shared-string order: note first, then new bullet,
# matching how the strings were appended to sharedStrings.xml).
$ws.Range("A6").Value = "note: 1 Delivery Note or Receipt Note will be Controlled in sales and purchase"
$ws.Range("A5").Value = "Stock Journal Import For Make Stock Item Closing balance zero"

# Add the "Cost" column header next to "Customization", bold like A1.
$ws.Range("B1").Value = "Cost"
$ws.Range("B1").Font.Bold = $true

# The old cost figure (3000) moves from B5 down to B7, with an updated
# amount.
$ws.Range("B5").ClearContents()
$ws.Range("B7").Value = 4500
